$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
# Copy header style (bold, border, alignment) from an existing header cell (A1) to the newly added header cells
$ws.Range("H1").Value = "Avg_Experiment_Time"
$ws.Range("I1").Value = "Std_Total_Rounds"
$ws.Range("J1").Value = "Std_Expl_Cost"
$ws.Range("K1").Value = "Std_Expl_Eff"
$ws.Range("L1").Value = "Std_Round_Time"
$ws.Range("M1").Value = "Std_Agent_Step_Time"
$ws.Range("N1").Value = "Std_Experiment_Time"
$ws.Range("O1").Value = "Obs_Prob"

# Apply the same style as the other header cells (bold font, border, centered) to the newly created header cells I1:O1
$ws.Range("A1").Copy() | Out-Null
$ws.Range("I1:O1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Update data rows (rows 2-13) with recomputed statistics and new Std_* columns ---
# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 55.282
$ws.Range("D2").Value = 55.282
$ws.Range("E2").Value = 3.09500908
$ws.Range("F2").Value = 0.1989465
$ws.Range("G2").Value = 0.1989465
$ws.Range("H2").Value = 10.84994284
$ws.Range("I2").Value = 6.128097842424896
$ws.Range("J2").Value = 6.128097842424896
$ws.Range("K2").Value = 0.3465613092005847
$ws.Range("L2").Value = 0.02770936486736398
$ws.Range("M2").Value = 0.02770936486736398
$ws.Range("N2").Value = 0.7396395379633052
$ws.Range("O2").Value = 0.15

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 90.408
$ws.Range("D3").Value = 90.408
$ws.Range("E3").Value = 1.90116602
$ws.Range("F3").Value = 0.13269806
$ws.Range("G3").Value = 0.13269806
$ws.Range("H3").Value = 11.82647556
$ws.Range("I3").Value = 11.93173367692533
$ws.Range("J3").Value = 11.93173367692533
$ws.Range("K3").Value = 0.2460433212795098
$ws.Range("L3").Value = 0.01984652958856585
$ws.Range("M3").Value = 0.01984652958856585
$ws.Range("N3").Value = 1.211475526843596
$ws.Range("O3").Value = 0.85

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 31.058
$ws.Range("D4").Value = 62.094
$ws.Range("E4").Value = 2.81651776
$ws.Range("F4").Value = 0.25050596
$ws.Range("G4").Value = 0.1252528
$ws.Range("H4").Value = 3.77786416
$ws.Range("I4").Value = 5.730996828589124
$ws.Range("J4").Value = 11.46455564215833
$ws.Range("K4").Value = 0.5294053305615996
$ws.Range("L4").Value = 0.05558009970344142
$ws.Range("M4").Value = 0.02779003914915961
$ws.Range("N4").Value = 0.5679643988450378
$ws.Range("O4").Value = 0.15

# Row 5
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 48.712
$ws.Range("D5").Value = 95.236
$ws.Range("E5").Value = 1.82073864
$ws.Range("F5").Value = 0.1785273
$ws.Range("G5").Value = 0.0892638
$ws.Range("H5").Value = 4.2477316
$ws.Range("I5").Value = 8.477073643369522
$ws.Range("J5").Value = 15.26835122395985
$ws.Range("K5").Value = 0.2970759085780355
$ws.Range("L5").Value = 0.03290491233572342
$ws.Range("M5").Value = 0.01645246477444363
$ws.Range("N5").Value = 0.54541678584841
$ws.Range("O5").Value = 0.85

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 15.136
$ws.Range("D6").Value = 60.514
$ws.Range("E6").Value = 2.9572709
$ws.Range("F6").Value = 0.29980848
$ws.Range("G6").Value = 0.07495210000000001
$ws.Range("H6").Value = 1.09935954
$ws.Range("I6").Value = 3.746140833215506
$ws.Range("J6").Value = 14.98823199034472
$ws.Range("K6").Value = 0.7049223475481079
$ws.Range("L6").Value = 0.07687250283341467
$ws.Range("M6").Value = 0.01921827907725826
$ws.Range("N6").Value = 0.2682806637338131
$ws.Range("O6").Value = 0.15

# Row 7
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 25.678
$ws.Range("D7").Value = 95.56999999999999
$ws.Range("E7").Value = 1.82186828
$ws.Range("F7").Value = 0.24335592
$ws.Range("G7").Value = 0.06083898000000001
$ws.Range("H7").Value = 1.52148324
$ws.Range("I7").Value = 5.393746471795278
$ws.Range("J7").Value = 16.6816562520219
$ws.Range("K7").Value = 0.3165163245902268
$ws.Range("L7").Value = 0.05475961097873765
$ws.Range("M7").Value = 0.01368999501411791
$ws.Range("N7").Value = 0.3007329667043972
$ws.Range("O7").Value = 0.85

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 9.682
$ws.Range("D8").Value = 58
$ws.Range("E8").Value = 3.10125016
$ws.Range("F8").Value = 0.35805146
$ws.Range("G8").Value = 0.05967514
$ws.Range("H8").Value = 0.56847574
$ws.Range("I8").Value = 2.518245245848083
$ws.Range("J8").Value = 15.11333271159285
$ws.Range("K8").Value = 0.7586036283162397
$ws.Range("L8").Value = 0.1066347532057486
$ws.Range("M8").Value = 0.01777221345027806
$ws.Range("N8").Value = 0.2045253984385933
$ws.Range("O8").Value = 0.15

# Row 9
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 17.088
$ws.Range("D9").Value = 88.09999999999999
$ws.Range("E9").Value = 1.98654412
$ws.Range("F9").Value = 0.25504854
$ws.Range("G9").Value = 0.04250812
$ws.Range("H9").Value = 0.7108454000000001
$ws.Range("I9").Value = 4.165448105046043
$ws.Range("J9").Value = 16.44417287805278
$ws.Range("K9").Value = 0.3814642765166121
$ws.Range("L9").Value = 0.05324956946835844
$ws.Range("M9").Value = 0.008874801443803341
$ws.Range("N9").Value = 0.1726952920186617
$ws.Range("O9").Value = 0.85

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 7.244
$ws.Range("D10").Value = 57.808
$ws.Range("E10").Value = 3.1094637
$ws.Range("F10").Value = 0.31977384
$ws.Range("G10").Value = 0.03997164
$ws.Range("H10").Value = 0.2897661599999999
$ws.Range("I10").Value = 1.828428369010173
$ws.Range("J10").Value = 14.6044416336559
$ws.Range("K10").Value = 0.7860699164320875
$ws.Range("L10").Value = 0.09241374735316228
$ws.Range("M10").Value = 0.01155152151922934
$ws.Range("N10").Value = 0.1197893736126683
$ws.Range("O10").Value = 0.15

# Row 11
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 13.072
$ws.Range("D11").Value = 81.604
$ws.Range("E11").Value = 2.146199220000001
$ws.Range("F11").Value = 0.24522004
$ws.Range("G11").Value = 0.03065262
$ws.Range("H11").Value = 0.39434694
$ws.Range("I11").Value = 3.469710908357518
$ws.Range("J11").Value = 15.67643726650738
$ws.Range("K11").Value = 0.4069580778062713
$ws.Range("L11").Value = 0.0492322764137528
$ws.Range("M11").Value = 0.006154130122993946
$ws.Range("N11").Value = 0.113916739230854
$ws.Range("O11").Value = 0.85

# Row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 5.826
$ws.Range("D12").Value = 58.108
$ws.Range("E12").Value = 3.1327158
$ws.Range("F12").Value = 0.29287424
$ws.Range("G12").Value = 0.02928736
$ws.Range("H12").Value = 0.172486
$ws.Range("I12").Value = 1.632499344017152
$ws.Range("J12").Value = 16.28246213403527
$ws.Range("K12").Value = 0.8603027700272695
$ws.Range("L12").Value = 0.09142711396713296
$ws.Range("M12").Value = 0.009142848614736217
$ws.Range("N12").Value = 0.07780056739493579
$ws.Range("O12").Value = 0.15

# Row 13
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 10.87
$ws.Range("D13").Value = 75.97799999999999
$ws.Range("E13").Value = 2.3149498
$ws.Range("F13").Value = 0.22915482
$ws.Range("G13").Value = 0.0229155
$ws.Range("H13").Value = 0.24672734
$ws.Range("I13").Value = 3.462524823927278
$ws.Range("J13").Value = 15.43431086585364
$ws.Range("K13").Value = 0.4635590358491514
$ws.Range("L13").Value = 0.045648012239422
$ws.Range("M13").Value = 0.004565058451688726
$ws.Range("N13").Value = 0.08703436025648095
$ws.Range("O13").Value = 0.85
